# Scheduled-runner data refresh: push updated market-price-derived
# profit figures (currentAveragePrice* / LevePrice* / LeveProfit*
# columns H-N) into each job sheet, row by row.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1685.0256
$ws.Range("I15").Value = 1685.0256
$ws.Range("K15").Value = 5055.0768
$ws.Range("M15").Value = -4886.0768
$ws.Range("H17").Value = 3062.1667
$ws.Range("J17").Value = 3062.1667
$ws.Range("L17").Value = 9186.500100000001
$ws.Range("N17").Value = -9522.500100000001
$ws.Range("H46").Value = 3950.5
$ws.Range("J46").Value = 3950.5
$ws.Range("L46").Value = 11851.5
$ws.Range("N46").Value = -12089.5
$ws.Range("H60").Value = 3950.5
$ws.Range("J60").Value = 3950.5
$ws.Range("L60").Value = 11851.5
$ws.Range("N60").Value = -12819.5
$ws.Range("H74").Value = 7089.619
$ws.Range("I74").Value = 4760.154
$ws.Range("J74").Value = 10875
$ws.Range("K74").Value = 4760.154
$ws.Range("L74").Value = 10875
$ws.Range("M74").Value = -3824.154
$ws.Range("N74").Value = -12747
$ws.Range("H77").Value = 7089.619
$ws.Range("I77").Value = 4760.154
$ws.Range("J77").Value = 10875
$ws.Range("K77").Value = 23800.77
$ws.Range("L77").Value = 54375
$ws.Range("M77").Value = -19120.77
$ws.Range("N77").Value = -63735
$ws.Range("H92").Value = 4630770.5
$ws.Range("I92").Value = 905.7895
$ws.Range("K92").Value = 905.7895
$ws.Range("M92").Value = 342.2105
$ws.Range("H132").Value = 1821.7894
$ws.Range("I132").Value = 1835.0714
$ws.Range("J132").Value = 1784.6
$ws.Range("K132").Value = 5505.2142
$ws.Range("L132").Value = 5353.799999999999
$ws.Range("M132").Value = -2975.2142
$ws.Range("N132").Value = -10413.8
$ws.Range("H137").Value = 2099.3428
$ws.Range("I137").Value = 1662.55
$ws.Range("J137").Value = 2681.7334
$ws.Range("K137").Value = 4987.65
$ws.Range("L137").Value = 8045.2002
$ws.Range("M137").Value = -2437.65
$ws.Range("N137").Value = -13145.2002
$ws.Range("H138").Value = 2914.258
$ws.Range("I138").Value = 1923.3334
$ws.Range("J138").Value = 2964.644
$ws.Range("K138").Value = 5770.0002
$ws.Range("L138").Value = 8893.931999999999
$ws.Range("M138").Value = -630.0002000000004
$ws.Range("N138").Value = -19173.932

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16405193
$ws.Range("I32").Value = 17552978
$ws.Range("K32").Value = 17552978
$ws.Range("M32").Value = -17552691
$ws.Range("H45").Value = 4349.3125
$ws.Range("I45").Value = 2929.5454
$ws.Range("K45").Value = 2929.5454
$ws.Range("M45").Value = -2552.5454
$ws.Range("H64").Value = 33570.715
$ws.Range("J64").Value = 33570.715
$ws.Range("L64").Value = 33570.715
$ws.Range("N64").Value = -34066.715
$ws.Range("H67").Value = 33570.715
$ws.Range("J67").Value = 33570.715
$ws.Range("L67").Value = 33570.715
$ws.Range("N67").Value = -35286.715
$ws.Range("H74").Value = 1644.2916
$ws.Range("I74").Value = 1233.75
$ws.Range("J74").Value = 3697
$ws.Range("K74").Value = 1233.75
$ws.Range("L74").Value = 3697
$ws.Range("M74").Value = -359.75
$ws.Range("N74").Value = -5445
$ws.Range("H77").Value = 1644.2916
$ws.Range("I77").Value = 1233.75
$ws.Range("J77").Value = 3697
$ws.Range("K77").Value = 6168.75
$ws.Range("L77").Value = 18485
$ws.Range("M77").Value = -1800.75
$ws.Range("N77").Value = -27221
$ws.Range("H132").Value = 2721.2727
$ws.Range("I132").Value = 2386.9673
$ws.Range("K132").Value = 7160.901899999999
$ws.Range("M132").Value = -4630.901899999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7694822
$ws.Range("I94").Value = 1893.4762
$ws.Range("J94").Value = 40005120
$ws.Range("K94").Value = 1893.4762
$ws.Range("L94").Value = 40005120
$ws.Range("M94").Value = -1442.4762
$ws.Range("N94").Value = -40006022
$ws.Range("H105").Value = 4650.25
$ws.Range("I105").Value = 2845.5715
$ws.Range("K105").Value = 2845.5715
$ws.Range("M105").Value = -1098.5715
$ws.Range("H134").Value = 3164.879
$ws.Range("I134").Value = 2031.6957
$ws.Range("K134").Value = 6095.0871
$ws.Range("M134").Value = -3560.0871
$ws.Range("H135").Value = 69997.414
$ws.Range("J135").Value = 69997.414
$ws.Range("L135").Value = 69997.414
$ws.Range("N135").Value = -80137.414

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3277.5945
$ws.Range("I31").Value = 2429.0688
$ws.Range("J31").Value = 6353.5
$ws.Range("K31").Value = 2429.0688
$ws.Range("L31").Value = 6353.5
$ws.Range("M31").Value = -2134.0688
$ws.Range("N31").Value = -6943.5
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("H34").Value = 3277.5945
$ws.Range("I34").Value = 2429.0688
$ws.Range("J34").Value = 6353.5
$ws.Range("K34").Value = 2429.0688
$ws.Range("L34").Value = 6353.5
$ws.Range("M34").Value = -2227.0688
$ws.Range("N34").Value = -6757.5
$ws.Range("H58").Value = 2107.2222
$ws.Range("I58").Value = 1103.1
$ws.Range("J58").Value = 3362.375
$ws.Range("K58").Value = 1103.1
$ws.Range("L58").Value = 3362.375
$ws.Range("M58").Value = -900.0999999999999
$ws.Range("N58").Value = -3768.375
$ws.Range("H86").Value = 25393.6
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 38989.332
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 38989.332
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -41235.332
$ws.Range("H89").Value = 25393.6
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 38989.332
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 194946.66
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -206178.66
$ws.Range("H132").Value = 3114.6177
$ws.Range("I132").Value = 2934.2812
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 8802.8436
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -6272.8436
$ws.Range("N132").Value = -23060
$ws.Range("H134").Value = 4381.375
$ws.Range("I134").Value = 2572.1875
$ws.Range("K134").Value = 7716.5625
$ws.Range("M134").Value = -5181.5625
$ws.Range("H136").Value = 2107.2222
$ws.Range("I136").Value = 1103.1
$ws.Range("J136").Value = 3362.375
$ws.Range("K136").Value = 3309.3
$ws.Range("L136").Value = 10087.125
$ws.Range("M136").Value = -759.2999999999997
$ws.Range("N136").Value = -15187.125
$ws.Range("M32").ClearContents()  # no NQ price data -> LeveProfitNQ is blank

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 157.18182
$ws.Range("I14").Value = 157.18182
$ws.Range("K14").Value = 471.5454599999999
$ws.Range("M14").Value = -298.5454599999999
$ws.Range("H56").Value = 7658.4116
$ws.Range("I56").Value = 7658.4116
$ws.Range("K56").Value = 7658.4116
$ws.Range("M56").Value = -7128.4116
$ws.Range("H122").Value = 2189.5454
$ws.Range("I122").Value = 1819.8
$ws.Range("J122").Value = 2497.6667
$ws.Range("K122").Value = 16378.2
$ws.Range("L122").Value = 22479.0003
$ws.Range("M122").Value = -13928.2
$ws.Range("N122").Value = -27379.0003

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 10870210
$ws.Range("I97").Value = 571.3889
$ws.Range("K97").Value = 571.3889
$ws.Range("M97").Value = -75.38890000000004
$ws.Range("H113").Value = 7109.143
$ws.Range("I113").Value = 4447.5557
$ws.Range("K113").Value = 4447.5557
$ws.Range("M113").Value = -2277.5557
$ws.Range("H132").Value = 3930.0613
$ws.Range("I132").Value = 3097.3513
$ws.Range("K132").Value = 9292.053899999999
$ws.Range("M132").Value = -6762.053899999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4936.154
$ws.Range("I7").Value = 2778.75
$ws.Range("J7").Value = 5895
$ws.Range("K7").Value = 2778.75
$ws.Range("L7").Value = 5895
$ws.Range("M7").Value = -2666.75
$ws.Range("N7").Value = -6119
$ws.Range("H32").Value = 14540
$ws.Range("I32").Value = 5233.3335
$ws.Range("K32").Value = 5233.3335
$ws.Range("M32").Value = -4916.3335
$ws.Range("H40").Value = 11169.685
$ws.Range("I40").Value = 14313.223
$ws.Range("K40").Value = 14313.223
$ws.Range("M40").Value = -14177.223
$ws.Range("H61").Value = 4694.727
$ws.Range("I61").Value = 2940.2
$ws.Range("K61").Value = 2940.2
$ws.Range("M61").Value = -2738.2
$ws.Range("H82").Value = 1253.3846
$ws.Range("J82").Value = 2735.3333
$ws.Range("L82").Value = 2735.3333
$ws.Range("N82").Value = -3457.3333
$ws.Range("H85").Value = 1253.3846
$ws.Range("J85").Value = 2735.3333
$ws.Range("L85").Value = 2735.3333
$ws.Range("N85").Value = -5231.3333
$ws.Range("H113").Value = 4694.727
$ws.Range("I113").Value = 2940.2
$ws.Range("K113").Value = 2940.2
$ws.Range("M113").Value = -770.1999999999998
$ws.Range("H126").Value = 4936.154
$ws.Range("I126").Value = 2778.75
$ws.Range("J126").Value = 5895
$ws.Range("K126").Value = 8336.25
$ws.Range("L126").Value = 17685
$ws.Range("M126").Value = -5866.25
$ws.Range("N126").Value = -22625
$ws.Range("H138").Value = 60473.332
$ws.Range("J138").Value = 60473.332
$ws.Range("L138").Value = 60473.332
$ws.Range("N138").Value = -70753.33199999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1022.61536
$ws.Range("I113").Value = 754
$ws.Range("K113").Value = 2262
$ws.Range("M113").Value = -92
$ws.Range("H126").Value = 1625.0476
$ws.Range("I126").Value = 1521.5625
$ws.Range("K126").Value = 4564.6875
$ws.Range("M126").Value = -2094.6875
$ws.Range("H132").Value = 3100.4
$ws.Range("I132").Value = 2556
$ws.Range("K132").Value = 7668
$ws.Range("M132").Value = -5138
